$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text does not look like a plain number (commas/dots
# patterns, percentages with padding spaces, URLs, names, subscript glyphs).
# These can be written directly with .Value and Excel keeps them as text,
# matching the original inlineStr cell type with no style change. ---

$ws.Range("D2").Value = "42.887.92"
$ws.Range("E2").Value = "  -1.23%  "
$ws.Range("D3").Value = "2.570.87"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("E13").Value = "  +6.49%  "
$ws.Range("D14").Value = "2.569.41"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "42.952.31"
$ws.Range("E17").Value = "  -0.99%  "

# Row 18 and row 19 swap their contents (Shiba Inu <-> Internet Computer)
$ws.Range("B18").Value = "InternetComputer(DFINITY)"
$ws.Range("C18").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E18").Value = "  +4.28%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0997"
$ws.Range("E19").Value = "  +2.84%  "

$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("E23").Value = "  +2.42%  "
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("E29").Value = "  -5.15%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("E31").Value = "  +3.53%  "
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  +13.26%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("E43").Value = "  +26.80%  "
$ws.Range("D44").Value = "2.072.77"
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("E48").Value = "  +12.95%  "
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "2.818.34"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("E51").Value = "  +1.67%  "

# --- Price cells in column D whose new text parses as a plain number
# (e.g. "97.11") would be auto-converted to a numeric cell by Excel's
# normal type inference. The source file keeps these as literal text
# (t="inlineStr"), so force text entry by flipping the cell to the "@"
# (Text) number format just long enough to type the value, then putting
# the cell's format back to the default "Normal" style so no visible
# formatting change remains. ---
$numericTextCells = [ordered]@{
    "D6"  = "97.11"
    "D9"  = "0.550"
    "D10" = "36.73"
    "D12" = "7.74"
    "D15" = "0.887"
    "D16" = "14.36"
    "D18" = "12.95"
    "D20" = "6.64"
    "D22" = "254.64"
    "D24" = "2.13"
    "D25" = "28.83"
    "D28" = "37.52"
    "D30" = "6.04"
    "D31" = "155.50"
    "D32" = "2.18"
    "D34" = "3.39"
    "D35" = "0.0808"
    "D36" = "18.36"
    "D39" = "23.51"
    "D40" = "3.44"
    "D43" = "2.04"
    "D45" = "1.00"
    "D46" = "9.27"
    "D47" = "85.48"
    "D48" = "76.97"
    "D49" = "106.47"
    "D51" = "0.192"
}

foreach ($addr in $numericTextCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $numericTextCells[$addr]
    $rng.Style = "Normal"
}
